$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.259705424308777
$ws.Range("B1").Value = 2.889725208282471
$ws.Range("C1").Value = 4.130145072937012
$ws.Range("D1").Value = 0.4835052788257599
$ws.Range("E1").Value = 0.6254622340202332
